# Overview.xlsx update
# The "DOM_GSEC" benchmark reference ("Breimann et al, 2023c") is bumped
# to the 2024 publication year in the "Reference" column (column I),
# which appears on both the DOM_GSEC and DOM_GSEC_PU rows (14 and 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = "Breimann et al, 2024c"
$ws.Range("I15").Value = "Breimann et al, 2024c"

# Leave the sheet with the cursor parked just below the used range,
# matching the saved selection state in the edited workbook.
$ws.Range("I17").Select()
